$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate the K column (strikeouts) values for the affected rows
$ws.Range("G2").Value = 3
$ws.Range("G5").Value = 0
$ws.Range("G9").Value = 2
$ws.Range("G10").Value = 3
